$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so numeric-looking strings
# (e.g. "1.00", "9.00", "561.57") are not coerced into numbers,
# matching the original inlineStr cell content.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.452.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.337.33"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.57"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.59"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.339.80"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.78%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.434"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.911.69"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.20%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.75"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.31%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.451.29"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.322.36"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.73%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.98%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.38"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.32"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.94"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.85%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0945"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.60"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +7.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.98"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.98%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.84"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.31"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.71"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +8.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.98"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.34%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +12.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.81"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0736"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.785.57"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.57%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.26"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.44"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.742"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.33%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.88%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.00"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.84%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "RenzoRestakedETH"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.377.28"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.104"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.84%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "287.24"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.09%  "
